$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("B1").Value = "Active cases"

# Cluster name / active cases data (rows 2-37)
$ws.Range("A2").Value = '3398 BlueCross Elly Kay Mordialloc'
$ws.Range("B2").Value = 31
$ws.Range("A3").Value = '3601 Baptcare Westhaven community'
$ws.Range("B3").Value = 13
$ws.Range("A4").Value = '3647 Aurrum Aged Care Reservoir'
$ws.Range("B4").Value = 11
$ws.Range("A5").Value = '3653 Fronditha Thalpori St Albans Aged Care'
$ws.Range("B5").Value = 21
$ws.Range("A6").Value = '3749 Rosebrook - McKenzie Aged CareRosebud'
$ws.Range("B6").Value = 11
$ws.Range("A7").Value = '3975 Aurrum Aged Care Brunswick West'
$ws.Range("B7").Value = 11
$ws.Range("A8").Value = '4257 BlueCross The Gables Camberwell'
$ws.Range("B8").Value = 16
$ws.Range("A9").Value = '4295 Hope Aged Care Sunshine West'
$ws.Range("B9").Value = 15
$ws.Range("A10").Value = '4314 Estia Health Ardeer'
$ws.Range("B10").Value = 16
$ws.Range("A11").Value = '44095 Myrniong Primary School Myrniong'
$ws.Range("B11").Value = 13
$ws.Range("A12").Value = '44404 Castlemaine North Primary SchoolCastlemaine'
$ws.Range("B12").Value = 61
$ws.Range("A13").Value = '44593 Torquay P-6 College Torquay'
$ws.Range("B13").Value = 19
$ws.Range("A14").Value = '44622 Grey Street Primary School Traralgon'
$ws.Range("B14").Value = 13
$ws.Range("A15").Value = '44631 Mount Evelyn Primary School'
$ws.Range("B15").Value = 23
$ws.Range("A16").Value = '44642 Irymple South Primary School Irymple South'
$ws.Range("B16").Value = 17
$ws.Range("A17").Value = '4479 Whittlesea Lodge Whittlesea'
$ws.Range("B17").Value = 15
$ws.Range("A18").Value = '44893 Greenhills Primary School Greensborough'
$ws.Range("B18").Value = 12
$ws.Range("A19").Value = '45168 Ranfurly Primary School Mildura'
$ws.Range("B19").Value = 28
$ws.Range("A20").Value = '45275 Lalor Gardens Primary School Lalor'
$ws.Range("B20").Value = 11
$ws.Range("A21").Value = '46320 St Mary''s Coptic Orthodox CollegeCoolaroo'
$ws.Range("B21").Value = 10
$ws.Range("A22").Value = '52390 Our Lady of the Way Catholic PrimarySchool Wallan'
$ws.Range("B22").Value = 49
$ws.Range("A23").Value = '52777 Mirripoa Primary School Mount DuneedSchool Camp'
$ws.Range("B23").Value = 31
$ws.Range("A24").Value = 'Confirmed Omicron Sircuit Bar Fitzroy'
$ws.Range("B24").Value = 19
$ws.Range("A25").Value = 'Confirmed Omicron Variant The Peel HotelCollingwood'
$ws.Range("B25").Value = 14
$ws.Range("A26").Value = 'Green Gables Lodge Warburton'
$ws.Range("B26").Value = 26
$ws.Range("A27").Value = 'Greendale Hotel Greendale'
$ws.Range("B27").Value = 15
$ws.Range("A28").Value = 'JBS Australia Brooklyn'
$ws.Range("B28").Value = 16
$ws.Range("A29").Value = 'Lockington Consolidated SchoolLockington'
$ws.Range("B29").Value = 16
$ws.Range("A30").Value = 'PGL Camp Rumbug Foster North'
$ws.Range("B30").Value = 11
$ws.Range("A31").Value = 'Social Gathering 11 Dec Windsor'
$ws.Range("B31").Value = 13
$ws.Range("A32").Value = 'St Brigid''s Parish Primary School Mordialloc'
$ws.Range("B32").Value = 13
$ws.Range("A33").Value = 'St Pauls Cathedral'
$ws.Range("B33").Value = 27
$ws.Range("A34").Value = 'St Vincents Hospital Melbourne EmergencyDepartment Fitzroy'
$ws.Range("B34").Value = 15
$ws.Range("A35").Value = 'St. Vincent''s Hospital Melbourne Fitzroy'
$ws.Range("B35").Value = 17
$ws.Range("A36").Value = 'StarTrack- Melbourne Tullamarine'
$ws.Range("B36").Value = 10
$ws.Range("A37").Value = 'Warbuton Lodge Warbuton'
$ws.Range("B37").Value = 14

# Remove now-unused trailing rows (old sheet had 49 rows of data)
$ws.Rows("38:49").Delete()
